$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("gens")
$ws1.Cells.Item(10, 3).Value = 30.57085
$ws1.Cells.Item(11, 3).Value = 100
$ws1.Cells.Item(12, 3).Value = 100
$ws1.Cells.Item(32, 3).Value = 155
$ws1.Cells.Item(33, 3).Value = 155
$ws1.Cells.Item(34, 3).Value = 350

$ws2 = $wb.Worksheets.Item("lines")
$ws2.Cells.Item(2, 3).Value = 28.233678
$ws2.Cells.Item(2, 4).Value = 0.06610408500000001
$ws2.Cells.Item(3, 3).Value = -98.478176
$ws2.Cells.Item(3, 4).Value = 0.62524317
$ws2.Cells.Item(4, 3).Value = -37.755503
$ws2.Cells.Item(4, 4).Value = 0.2915129
$ws2.Cells.Item(5, 3).Value = -49.481842
$ws2.Cells.Item(5, 4).Value = 0.33501137
$ws2.Cells.Item(6, 3).Value = -19.28448
$ws2.Cells.Item(6, 4).Value = 0.40009643
$ws2.Cells.Item(7, 3).Value = 17.208918
$ws2.Cells.Item(7, 4).Value = 0.050212304
$ws2.Cells.Item(8, 3).Value = -295.68709
$ws2.Cells.Item(8, 4).Value = 0.21346799
$ws2.Cells.Item(9, 3).Value = -123.48184
$ws2.Cells.Item(9, 4).Value = 0.27434002
$ws2.Cells.Item(10, 3).Value = -108.7555
$ws2.Cells.Item(10, 4).Value = 0.30180159
$ws2.Cells.Item(11, 3).Value = -155.28448
$ws2.Cells.Item(11, 4).Value = 0.12711397
$ws2.Cells.Item(12, 3).Value = 105.57085
$ws2.Cells.Item(13, 3).Value = -50.748914
$ws2.Cells.Item(13, 4).Value = 0.041070496
$ws2.Cells.Item(13, 5).Value = 0
$ws2.Cells.Item(14, 3).Value = -14.680236
$ws2.Cells.Item(14, 4).Value = -0.041070496
$ws2.Cells.Item(15, 3).Value = -159.11116
$ws2.Cells.Item(15, 4).Value = 0.12235833
$ws2.Cells.Item(16, 3).Value = -172.91068
$ws2.Cells.Item(16, 4).Value = 0.15557659
$ws2.Cells.Item(17, 3).Value = -229.96035
$ws2.Cells.Item(17, 4).Value = 0.20449932
$ws2.Cells.Item(18, 3).Value = -243.75987
$ws2.Cells.Item(18, 4).Value = 0.23771759
$ws2.Cells.Item(19, 3).Value = -233.17131
$ws2.Cells.Item(19, 4).Value = 0.092496465
$ws2.Cells.Item(20, 3).Value = -155.90019
$ws2.Cells.Item(20, 4).Value = 0.082494416
$ws2.Cells.Item(21, 3).Value = -209.02215
$ws2.Cells.Item(21, 4).Value = 0.059278198
$ws2.Cells.Item(22, 3).Value = -207.6484
$ws2.Cells.Item(22, 4).Value = 0.33436978
$ws2.Cells.Item(23, 3).Value = -116.19346
$ws2.Cells.Item(23, 4).Value = 0.27509158
$ws2.Cells.Item(24, 3).Value = -349.90019
$ws2.Cells.Item(24, 4).Value = 0.11588501
$ws2.Cells.Item(25, 3).Value = 100.23817
$ws2.Cells.Item(25, 4).Value = 0.025335212
$ws2.Cells.Item(26, 3).Value = -278.96263
$ws2.Cells.Item(26, 4).Value = 0.025748986
$ws2.Cells.Item(27, 3).Value = -278.96263
$ws2.Cells.Item(27, 4).Value = 0.025748986
$ws2.Cells.Item(28, 3).Value = 295.68709
$ws2.Cells.Item(28, 4).Value = -0.13214685
$ws2.Cells.Item(29, 3).Value = -433.7336
$ws2.Cells.Item(29, 4).Value = 0.008345129099999999
$ws2.Cells.Item(29, 5).Value = 0
$ws2.Cells.Item(30, 3).Value = 239.07158
$ws2.Cells.Item(30, 4).Value = 0.072070335
$ws2.Cells.Item(31, 3).Value = -292.13739
$ws2.Cells.Item(31, 4).Value = 0.0051353748
$ws2.Cells.Item(32, 3).Value = -141.59621
$ws2.Cells.Item(32, 4).Value = -0.0048138281
$ws2.Cells.Item(33, 3).Value = -0.23926241
$ws2.Cells.Item(33, 4).Value = -0.01306673
$ws2.Cells.Item(34, 3).Value = -0.23926241
$ws2.Cells.Item(34, 4).Value = -0.01306673
$ws2.Cells.Item(35, 3).Value = 29.035789
$ws2.Cells.Item(35, 4).Value = 0.062669857
$ws2.Cells.Item(36, 3).Value = 29.035789
$ws2.Cells.Item(36, 4).Value = 0.062669857
$ws2.Cells.Item(37, 3).Value = -34.964211
$ws2.Cells.Item(37, 4).Value = 0.034468421
$ws2.Cells.Item(38, 3).Value = -34.964211
$ws2.Cells.Item(38, 4).Value = 0.034468421
$ws2.Cells.Item(39, 3).Value = -158.40379
$ws2.Cells.Item(39, 4).Value = 0.0031175268
$ws2.Cells.Item(40, 3).Value = 266.22972
$ws2.Cells.Item(40, 4).Value = 0.087683185

$ws3 = $wb.Worksheets.Item("bus")
$ws3.Cells.Item(2, 2).Value = 99.36561500000001
$ws3.Cells.Item(2, 3).Value = 14.084416
$ws3.Cells.Item(3, 2).Value = 99.431719
$ws3.Cells.Item(3, 3).Value = 13.689145
$ws3.Cells.Item(4, 2).Value = 99.990858
$ws3.Cells.Item(4, 3).Value = 34.863311
$ws3.Cells.Item(5, 2).Value = 99.76673
$ws3.Cells.Item(5, 3).Value = 19.973338
$ws3.Cells.Item(6, 2).Value = 99.657128
$ws3.Cells.Item(6, 3).Value = 17.293634
$ws3.Cells.Item(7, 2).Value = 99.831816
$ws3.Cells.Item(7, 3).Value = 17.391765
$ws3.Cells.Item(8, 3).Value = 30.881701
$ws3.Cells.Item(9, 3).Value = 24.441879
$ws3.Cells.Item(10, 2).Value = 100.04107
$ws3.Cells.Item(10, 3).Value = 32.81545
$ws3.Cells.Item(11, 2).Value = 99.95893
$ws3.Cells.Item(11, 3).Value = 26.864118
$ws3.Cells.Item(12, 2).Value = 100.16343
$ws3.Cells.Item(12, 3).Value = 46.180787
$ws3.Cells.Item(13, 2).Value = 100.19665
$ws3.Cells.Item(13, 3).Value = 47.339947
$ws3.Cells.Item(14, 2).Value = 100.25593
$ws3.Cells.Item(14, 3).Value = 57.37301
$ws3.Cells.Item(15, 2).Value = 100.24592
$ws3.Cells.Item(15, 3).Value = 52.728595
$ws3.Cells.Item(16, 2).Value = 100.33647
$ws3.Cells.Item(16, 3).Value = 75.076756
$ws3.Cells.Item(17, 2).Value = 100.36181
$ws3.Cells.Item(17, 3).Value = 73.37270700000001
$ws3.Cells.Item(18, 2).Value = 100.37015
$ws3.Cells.Item(18, 3).Value = 84.64978000000001
$ws3.Cells.Item(19, 2).Value = 100.37529
$ws3.Cells.Item(19, 3).Value = 88.739704
$ws3.Cells.Item(20, 2).Value = 100.43388
$ws3.Cells.Item(20, 3).Value = 67.874061
$ws3.Cells.Item(21, 2).Value = 100.49655
$ws3.Cells.Item(21, 3).Value = 66.71262900000001
$ws3.Cells.Item(22, 2).Value = 100.36222
$ws3.Cells.Item(22, 3).Value = 88.745925
$ws3.Cells.Item(23, 2).Value = 100.36534
$ws3.Cells.Item(23, 3).Value = 99.517383
$ws3.Cells.Item(24, 2).Value = 100.53102
$ws3.Cells.Item(24, 3).Value = 67.481842
$ws3.Cells.Item(25, 2).Value = 100.20433
$ws3.Cells.Item(25, 3).Value = 59.701027
$ws3.Cells.Item(26, 2).Value = 100.6187
$ws3.Cells.Item(26, 3).Value = 65.085774
